# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D) values, preserving original text representation ---
# (NumberFormat is temporarily forced to Text so Excel does not auto-convert
#  strings such as "1.000" or "310.60" into numeric values, then the cell
#  style is restored to its original so no formatting changes are introduced.)
$origStyle = $ws.Range('D2').Style
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '27.649.96'
$ws.Range('D2').Style = $origStyle
$origStyle = $ws.Range('D3').Style
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.895.68'
$ws.Range('D3').Style = $origStyle
$origStyle = $ws.Range('D4').Style
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = $origStyle
$origStyle = $ws.Range('D5').Style
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '310.60'
$ws.Range('D5').Style = $origStyle
$origStyle = $ws.Range('D6').Style
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = $origStyle
$origStyle = $ws.Range('D7').Style
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5248'
$ws.Range('D7').Style = $origStyle
$origStyle = $ws.Range('D8').Style
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3780'
$ws.Range('D8').Style = $origStyle
$origStyle = $ws.Range('D9').Style
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.07225'
$ws.Range('D9').Style = $origStyle
$origStyle = $ws.Range('D10').Style
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '21.04'
$ws.Range('D10').Style = $origStyle
$origStyle = $ws.Range('D11').Style
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.8996'
$ws.Range('D11').Style = $origStyle
$origStyle = $ws.Range('D12').Style
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.896.25'
$ws.Range('D12').Style = $origStyle
$origStyle = $ws.Range('D13').Style
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.07632'
$ws.Range('D13').Style = $origStyle
$origStyle = $ws.Range('D14').Style
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.421'
$ws.Range('D14').Style = $origStyle
$origStyle = $ws.Range('D15').Style
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '91.66'
$ws.Range('D15').Style = $origStyle
$origStyle = $ws.Range('D16').Style
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = $origStyle
$origStyle = $ws.Range('D17').Style
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.000008672'
$ws.Range('D17').Style = $origStyle
$origStyle = $ws.Range('D18').Style
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '14.25'
$ws.Range('D18').Style = $origStyle
$origStyle = $ws.Range('D20').Style
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '27.686.56'
$ws.Range('D20').Style = $origStyle
$origStyle = $ws.Range('D21').Style
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.139'
$ws.Range('D21').Style = $origStyle
$origStyle = $ws.Range('D22').Style
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.126.25'
$ws.Range('D22').Style = $origStyle
$origStyle = $ws.Range('D23').Style
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.81'
$ws.Range('D23').Style = $origStyle
$origStyle = $ws.Range('D25').Style
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '153.15'
$ws.Range('D25').Style = $origStyle
$origStyle = $ws.Range('D26').Style
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.859'
$ws.Range('D26').Style = $origStyle
$origStyle = $ws.Range('D27').Style
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '18.23'
$ws.Range('D27').Style = $origStyle
$origStyle = $ws.Range('D28').Style
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.163'
$ws.Range('D28').Style = $origStyle
$origStyle = $ws.Range('D29').Style
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '113.93'
$ws.Range('D29').Style = $origStyle
$origStyle = $ws.Range('D30').Style
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.805'
$ws.Range('D30').Style = $origStyle
$origStyle = $ws.Range('D31').Style
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.816'
$ws.Range('D31').Style = $origStyle
$origStyle = $ws.Range('D32').Style
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09152'
$ws.Range('D32').Style = $origStyle
$origStyle = $ws.Range('D33').Style
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.05267'
$ws.Range('D33').Style = $origStyle
$origStyle = $ws.Range('D34').Style
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.164'
$ws.Range('D34').Style = $origStyle
$origStyle = $ws.Range('D36').Style
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7726'
$ws.Range('D36').Style = $origStyle
$origStyle = $ws.Range('D37').Style
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02079'
$ws.Range('D37').Style = $origStyle
$origStyle = $ws.Range('D38').Style
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.565'
$ws.Range('D38').Style = $origStyle
$origStyle = $ws.Range('D40').Style
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.091'
$ws.Range('D40').Style = $origStyle
$origStyle = $ws.Range('D41').Style
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5555'
$ws.Range('D41').Style = $origStyle
$origStyle = $ws.Range('D42').Style
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '6.682'
$ws.Range('D42').Style = $origStyle
$origStyle = $ws.Range('D43').Style
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '117.02'
$ws.Range('D43').Style = $origStyle
$origStyle = $ws.Range('D44').Style
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '8.724'
$ws.Range('D44').Style = $origStyle
$origStyle = $ws.Range('D45').Style
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1511'
$ws.Range('D45').Style = $origStyle
$origStyle = $ws.Range('D46').Style
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.4798'
$ws.Range('D46').Style = $origStyle
$origStyle = $ws.Range('D47').Style
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '10.38'
$ws.Range('D47').Style = $origStyle
$origStyle = $ws.Range('D48').Style
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.9996'
$ws.Range('D48').Style = $origStyle
$origStyle = $ws.Range('D49').Style
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '66.22'
$ws.Range('D49').Style = $origStyle
$origStyle = $ws.Range('D50').Style
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.577'
$ws.Range('D50').Style = $origStyle

# --- Update "Volume(1h)" column (E) percentage text values ---
$ws.Range('E2').Value = '  -0.57%  '
$ws.Range('E3').Value = '  -0.19%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  +6.51%  '
$ws.Range('E8').Value = '  -1.01%  '
$ws.Range('E9').Value = '  -1.31%  '
$ws.Range('E10').Value = '  -0.54%  '
$ws.Range('E11').Value = '  -0.93%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('E18').Value = '  -1.84%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('E22').Value = '  +0.05%  '
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('E25').Value = '  -0.76%  '
$ws.Range('E26').Value = '  -0.20%  '
$ws.Range('E27').Value = '  -0.78%  '
$ws.Range('E28').Value = '  -2.85%  '
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  -2.02%  '
$ws.Range('E31').Value = '  +3.80%  '
$ws.Range('E32').Value = '  +2.38%  '
$ws.Range('E33').Value = '  -0.54%  '
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('E35').Value = '  -1.44%  '
$ws.Range('E36').Value = '  +0.13%  '
$ws.Range('E37').Value = '  +0.88%  '
$ws.Range('E38').Value = '  -0.39%  '
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('E41').Value = '  +0.39%  '
$ws.Range('E42').Value = '  -4.38%  '
$ws.Range('E43').Value = '  +5.14%  '
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('E45').Value = '  -0.63%  '
$ws.Range('E46').Value = '  +0.09%  '
$ws.Range('E47').Value = '  -2.71%  '
$ws.Range('E48').Value = '  -0.03%  '
$ws.Range('E49').Value = '  -1.81%  '
$ws.Range('E50').Value = '  -3.94%  '
$ws.Range('E51').Value = '  -0.12%  '
